$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-36): "Förändrad" date increases from 45668 to 45669 (one day later)
for ($r = 2; $r -le 36; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45668) {
        $cell.Value2 = 45669
    }
}

# Swap row 35 and row 36 values for column A (Beteckning) and column G (Area ha)
$a35 = $ws.Cells.Item(35, 1).Value2
$a36 = $ws.Cells.Item(36, 1).Value2
$ws.Cells.Item(35, 1).Value2 = $a36
$ws.Cells.Item(36, 1).Value2 = $a35

$g35 = $ws.Cells.Item(35, 7).Value2
$g36 = $ws.Cells.Item(36, 7).Value2
$ws.Cells.Item(35, 7).Value2 = $g36
$ws.Cells.Item(36, 7).Value2 = $g35
